$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '39.789.17'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.225.62'
$ws.Range('E3').Value = '  -5.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '294.13'
$ws.Range('E5').Value = '  -5.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '84.34'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.517'
$ws.Range('E7').Value = '  -2.49%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.467'
$ws.Range('E9').Value = '  -3.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0788'
$ws.Range('E10').Value = '  -2.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '29.89'
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.71'
$ws.Range('E12').Value = '  -9.03%  '
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.569.27'
$ws.Range('E14').Value = '  -5.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.32'
$ws.Range('E15').Value = '  -1.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.12'
$ws.Range('E16').Value = '  -4.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.225.19'
$ws.Range('E17').Value = '  -6.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.723'
$ws.Range('E18').Value = '  -4.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '39.708.14'
$ws.Range('E19').Value = '  -0.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0884'
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.77'
$ws.Range('E21').Value = '  -5.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.35'
$ws.Range('E22').Value = '  -4.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.52'
$ws.Range('E23').Value = '  -1.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '232.45'
$ws.Range('E24').Value = '  -1.07%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.41'
$ws.Range('E26').Value = '  -5.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.82'
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.91'
$ws.Range('E28').Value = '  -3.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').Value = '  +2.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.21'
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.53'
$ws.Range('E31').Value = '  -5.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '150.05'
$ws.Range('E32').Value = '  -2.42%  '
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.84'
$ws.Range('E34').Value = '  -5.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.37'
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0704'
$ws.Range('E36').Value = '  -1.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.07'
$ws.Range('E37').Value = '  +3.36%  '
$ws.Range('E38').Value = '  -2.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0979'
$ws.Range('E39').Value = '  -0.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.67'
$ws.Range('E40').Value = '  -4.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.66'
$ws.Range('E41').Value = '  -3.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.70'
$ws.Range('E42').Value = '  -4.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.946.73'
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('E44').Value = '  -3.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0267'
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.37'
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.29'
$ws.Range('E47').Value = '  -7.29%  '
$ws.Range('E48').Value = '  -3.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.434.63'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '70.86'
$ws.Range('E50').Value = '  +0.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '88.98'
$ws.Range('E51').Value = '  -4.44%  '
